# Fix the car title text on "findNewCarTest" (D3 was truncated to "Ho Cars",
# should read "Honda Cars"), then make "findNewCarTest" the active sheet with
# H6 selected (it was previously "carNameAndPrice" with H8 selected there).

$wb = $excel.ActiveWorkbook

$wsFindNewCar = $wb.Worksheets.Item("findNewCarTest")

# Correct the truncated car title.
$wsFindNewCar.Range("D3").Value = "Honda Cars"

# Switch the active tab / selection to findNewCarTest!H6.
$wsFindNewCar.Activate()
$wsFindNewCar.Range("H6").Select()
